$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..61 down to 9..62.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data entry.
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(8, 3).Value = "Maule"
$ws.Cells.Item(8, 4).Value = (Get-Date -Year 2022 -Month 9 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = 100112040
$ws.Cells.Item(8, 7).Value = "Cilantro"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 200
$ws.Cells.Item(8, 11).Value = 7500
$ws.Cells.Item(8, 12).Value = 7500
$ws.Cells.Item(8, 13).Value = 7500
$ws.Cells.Item(8, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(8, 15).Value = "Región del Maule"
$ws.Cells.Item(8, 16).Value = 208
$ws.Cells.Item(8, 17).Value = 36
$ws.Cells.Item(8, 18).Value = "Hortaliza"
